$wb = $excel.ActiveWorkbook
$s1 = $wb.Worksheets.Item(1)

# Insert a new row at row 11 (shifts old rows 11-15 down to 12-16),
# making room for the new second "Contact" entry (Bob Milius).
$s1.Rows.Item(11).Insert()

# Copy formatting from the row above into the freshly inserted row so it
# matches the rest of the data rows (border + vertical/wrap alignment),
# rather than getting a generic blank style.
$s1.Range("A10:B10").Copy()
$s1.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Metadata sheet updates ---

# Version: 0.1.6 -> 0.1.7
$s1.Range("B3").Value = "0.1.7"

# Status: active -> draft
$s1.Range("B6").Value = "draft"

# Date refreshed
$s1.Range("B8").Value = "2024-08-27T12:23:18-05:00"

# Contact (row 10): display text updated to the publisher org + website
$s1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# New row 11: second Contact entry (individual contact)
$s1.Range("A11").Value = "Contact"
$s1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Row 12 (previously row 11): now Jurisdiction / blank
$s1.Range("A12").Value = "Jurisdiction"
$s1.Range("B12").Value = ""

# Row 13 (previously row 12): Description
$s1.Range("A13").Value = "Description"
$s1.Range("B13").Value = "Platelet morphology panel - Blood (58406-0)"

# Row 14 (previously row 13): Purpose, blank - unchanged
$s1.Range("A14").Value = "Purpose"
$s1.Range("B14").Value = ""

# Row 15 (previously row 14): Copyright, blank - unchanged
$s1.Range("A15").Value = "Copyright"
$s1.Range("B15").Value = ""

# Row 16 (previously row 15): Immutable / BooleanType[null] - unchanged
$s1.Range("A16").Value = "Immutable"
$s1.Range("B16").Value = "BooleanType[null]"
